# Weekly refresh of the "Hortaliza, Terminal La Palmera de La Serena - Perejil"
# data table: a new week's record is inserted at the top of the data block
# (row 12), pushing all subsequent rows down by one, and the sheet's
# dimension grows from A1:R101 to A1:R102 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 12 (row 1 is the header, rows 2-11 stay
# fixed; this shifts the former rows 12-101 down to 13-102).
$ws.Rows.Item(12).Insert()

# The newly inserted row has no values yet, only the row-above's formatting
# carried down by Insert(). Seed it with the same repeating field set every
# other data row in this table uses (Mercado ID, Mercado, Región, Codreg,
# Categoría ID, Categoría, Variedad, Calidad, Unidad de comercialización,
# Origen, Kg o Unidades, Clasificación) by duplicating the row right below
# it (the row that used to be row 12 before the insert).
$ws.Rows.Item(13).Copy()
$ws.Rows.Item(12).PasteSpecial()

# Now overwrite the week-specific figures for the new row 12.
$ws.Cells.Item(12, 4).Value() = 44537   # Fecha
$ws.Cells.Item(12, 10).Value() = 3100   # Volumen
$ws.Cells.Item(12, 11).Value() = 1500   # Precio mínimo
$ws.Cells.Item(12, 12).Value() = 2000   # Precio máximo
$ws.Cells.Item(12, 13).Value() = 1750   # Precio promedio ponderado
$ws.Cells.Item(12, 16).Value() = 1167   # Precio $/Kg
